# Update the weekly schedule from the week of 20/05/2024 to the week of
# 08/06/2024 across every "Table N AM/PM" worksheet, and refresh a handful
# of booking counts that changed alongside the new dates.

$wb = $excel.ActiveWorkbook

# Old date (dd/mm/yyyy) -> new date (dd/mm/yyyy), row by row (rows 2-11).
$newDates = @(
    "08/06/2024",
    "09/06/2024",
    "10/06/2024",
    "11/06/2024",
    "12/06/2024",
    "13/06/2024",
    "14/06/2024",
    "15/06/2024",
    "16/06/2024",
    "17/06/2024"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    for ($r = 0; $r -lt $newDates.Length; $r++) {
        $row = $r + 2
        # Prefix with an apostrophe so the dd/mm/yyyy text is kept as text
        # instead of being auto-converted into a date serial number.
        $ws.Cells.Item($row, 1).Value = "'" + $newDates[$r]
    }
}

# Booking-count corrections that came in with the new week's data.
$wb.Worksheets.Item(1).Range("B2").Value = 2    # Table 1 AM
$wb.Worksheets.Item(4).Range("C2").Value = 4    # Table 2 PM
$wb.Worksheets.Item(5).Range("B2").Value = 4    # Table 3 AM
$wb.Worksheets.Item(5).Range("B4").Value = 4    # Table 3 AM
$wb.Worksheets.Item(5).Range("D9").Value = 4    # Table 3 AM
$wb.Worksheets.Item(6).Range("B2").Value = 4    # Table 3 PM
